$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 110, shifting existing rows 110:180 down to 111:181.
$ws.Rows.Item(110).Insert()

# Populate the newly inserted row 110 with the new record's data.
$ws.Range("A110").Value = 10
$ws.Range("B110").Value = "Vega Modelo de Temuco"
$ws.Range("C110").Value = "La Araucanía"
$ws.Range("D110").Value = 44603
$ws.Range("E110").Value = 9
$ws.Range("F110").Value = 100112043
$ws.Range("G110").Value = "Pepino dulce"
$ws.Range("H110").Value = "Cultivar IV Región"
$ws.Range("I110").Value = "Primera"
$ws.Range("J110").Value = 68
$ws.Range("K110").Value = 18000
$ws.Range("L110").Value = 18000
$ws.Range("M110").Value = 18000
$ws.Range("N110").Value = "$/bandeja 18 kilos"
$ws.Range("O110").Value = "Provincia de Limarí"
$ws.Range("P110").Value = 1000
$ws.Range("Q110").Value = 18
$ws.Range("R110").Value = "Hortaliza"
